$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number and date range) ---
$ws.Range("A8").Value = "Volume 31   Number  12"
$ws.Range("C9").Value = "Report Covering the Week  3/18/2024  Through  3/24/2024"

# --- Data table updates (rows 14-30) ---
# Row 14
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("F14").Value = 3
$ws.Range("N14").Value = -53.846153846153

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range("C15").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "***.*"
$ws.Range("C15").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("N15").Value = -63.636363636363

# Row 16
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -50
$ws.Range("G16").Value = 24
$ws.Range("H16").Value = 12.5
$ws.Range("I16").Value = 67
$ws.Range("J16").Value = 72
$ws.Range("K16").Value = -6.944444444444
$ws.Range("L16").Value = -17.283950617283
$ws.Range("M16").Value = -23.863636363636
$ws.Range("N16").Value = -87.773722627737

# Row 17
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 14
$ws.Range("E17").Value = -42.857142857142
$ws.Range("F17").Value = 36
$ws.Range("G17").Value = 39
$ws.Range("H17").Value = -7.692307692307
$ws.Range("I17").Value = 152
$ws.Range("J17").Value = 152
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = -0.653594771241
$ws.Range("M17").Value = 85.365853658536
$ws.Range("N17").Value = -26.570048309178

# Row 18
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 18
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = 28.571428571428
$ws.Range("I18").Value = 55
$ws.Range("J18").Value = 42
$ws.Range("K18").Value = 30.952380952381
$ws.Range("L18").Value = 7.843137254901
$ws.Range("M18").Value = 14.583333333333
$ws.Range("N18").Value = -88.565488565488

# Row 19
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 33.333333333333
$ws.Range("F19").Value = 54
$ws.Range("H19").Value = 35
$ws.Range("I19").Value = 156
$ws.Range("J19").Value = 119
$ws.Range("K19").Value = 31.092436974789
$ws.Range("L19").Value = 8.333333333333
$ws.Range("M19").Value = 188.888888888889
$ws.Range("N19").Value = -20.812182741116

# Row 20
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 14
$ws.Range("H20").Value = -36.363636363636
$ws.Range("I20").Value = 46
$ws.Range("J20").Value = 70
$ws.Range("K20").Value = -34.285714285714
$ws.Range("L20").Value = -31.343283582089
$ws.Range("M20").Value = 91.666666666666
$ws.Range("N20").Value = -81.526104417670

# Row 21
$ws.Range("C21").Value = 31
$ws.Range("E21").Value = -24.390243902439
$ws.Range("F21").Value = 153
$ws.Range("G21").Value = 142
$ws.Range("H21").Value = 7.746478873239
$ws.Range("I21").Value = 486
$ws.Range("J21").Value = 462
$ws.Range("K21").Value = 5.194805194805
$ws.Range("L21").Value = -3.187250996015
$ws.Range("M21").Value = 62
$ws.Range("N21").Value = -71.512309495896

# Row 22
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("D22").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("F22").Value = 1
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "0"
$ws.Range("D22").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "***.*"
$ws.Range("D22").Copy()
$ws.Range("H22").PasteSpecial(-4122)

# Row 23
$ws.Range("C23").Value = 2
$ws.Range("C23").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("D23").Value = 1
$ws.Range("H23").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 200
$ws.Range("I23").Value = 7
$ws.Range("J23").Value = 9
$ws.Range("K23").Value = -22.222222222222
$ws.Range("L23").Value = -12.5
$ws.Range("M23").Value = 40

# Row 24
$ws.Range("C24").Value = 10
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = -56.521739130434
$ws.Range("F24").Value = 73
$ws.Range("G24").Value = 68
$ws.Range("H24").Value = 7.352941176470
$ws.Range("I24").Value = 216
$ws.Range("J24").Value = 204
$ws.Range("K24").Value = 5.882352941176
$ws.Range("L24").Value = -27.027027027027
$ws.Range("M24").Value = -1.369863013698

# Row 25
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = -87.5
$ws.Range("F25").Value = 13
$ws.Range("G25").Value = 14
$ws.Range("H25").Value = -7.142857142857
$ws.Range("I25").Value = 48
$ws.Range("J25").Value = 58
$ws.Range("K25").Value = -17.241379310344
$ws.Range("L25").Value = -71.929824561403

# Row 26
$ws.Range("C26").Value = 16
$ws.Range("D26").Value = 21
$ws.Range("E26").Value = -23.809523809523
$ws.Range("F26").Value = 74
$ws.Range("G26").Value = 65
$ws.Range("H26").Value = 13.846153846153
$ws.Range("I26").Value = 257
$ws.Range("J26").Value = 253
$ws.Range("K26").Value = 1.581027667984
$ws.Range("L26").Value = 21.800947867298
$ws.Range("M26").Value = 0.390625

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("C27").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("C27").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("L27").Value = -10

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("M28").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "***.*"
$ws.Range("M28").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("I28").Value = 17
$ws.Range("K28").Value = -5.555555555555
$ws.Range("L28").Value = -10.526315789473

# Row 29
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "0"
$ws.Range("D29").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("F29").Value = 3
$ws.Range("H29").Value = 200
$ws.Range("M29").Value = -40
$ws.Range("N29").Value = -78.571428571428

# Row 30
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "0"
$ws.Range("D30").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("F30").Value = 3
$ws.Range("H30").Value = 200
$ws.Range("M30").Value = -40
$ws.Range("N30").Value = -76.923076923076

